$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell "Save" in H1, matching the bold/bordered/centered header style
# already used by the other header cells (e.g. G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Values for the new "Save" column, rows 2-19.
$saveValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    17 = 0
    18 = 1
    19 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
